# Generate Report for Handback
# Update the handoff/handback timestamps for the zh-cn and de-de reports
# to reflect the latest report generation run.

$wb = $excel.ActiveWorkbook

# zh-cn sheet: Correspond Handoff Datetime (E2) / Correspond Handback DateTime (H2)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-23 19:23:29"
$wsZhCn.Range("H2").Value = "2016-03-23 19:23:55"

# de-de sheet: Correspond Handoff Datetime (E2) / Correspond Handback DateTime (H2)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-23 19:23:35"
$wsDeDe.Range("H2").Value = "2016-03-23 19:24:01"
